$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 58, shifting rows 58:141 down to 59:142
$ws.Rows(58).Insert()

# Populate the new row 58 with the data for the inserted record
$ws.Cells.Item(58, 1).Value = 4
$ws.Cells.Item(58, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(58, 3).Value = "Los Lagos"
$ws.Cells.Item(58, 4).Value = 44482
$ws.Cells.Item(58, 5).Value = 10
$ws.Cells.Item(58, 6).Value = 100112017
$ws.Cells.Item(58, 7).Value = "Apio"
$ws.Cells.Item(58, 8).Value = "Americana (o)"
$ws.Cells.Item(58, 9).Value = "Primera"
$ws.Cells.Item(58, 10).Value = 20
$ws.Cells.Item(58, 11).Value = 11000
$ws.Cells.Item(58, 12).Value = 11000
$ws.Cells.Item(58, 13).Value = 11000
$ws.Cells.Item(58, 14).Value = "$/docena de matas"
$ws.Cells.Item(58, 15).Value = "Región de Coquimbo"
$ws.Cells.Item(58, 16).Value = 1833
$ws.Cells.Item(58, 17).Value = 6
$ws.Cells.Item(58, 18).Value = "Hortaliza"
